$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.309.61"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.875.33"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'0.7125"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'242.63"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.3112"
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("D9").Value = "'0.07768"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'25.08"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.08456"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").Value = "1.872.23"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "'5.206"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'0.7118"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "'91.33"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "29.310.89"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'0.000008322"
$ws.Range("E17").Value = "  +6.74%  "
$ws.Range("D18").Value = "'5.993"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").Value = "'242.58"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "2.122.39"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'7.823"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'0.1608"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "'162.79"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "'9.021"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "'1.514"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "'4.409"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").Value = "'4.335"
$ws.Range("E31").Value = "  +6.07%  "
$ws.Range("D32").Value = "'1.267"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").Value = "'0.05257"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "'1.923"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "'1.173"
$ws.Range("D36").Value = "'0.7466"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'0.01860"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'2.719"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "1.167.89"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "'6.371"
$ws.Range("E41").Value = "  +4.89%  "
$ws.Range("D42").Value = "'72.98"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "'0.8860"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "'106.57"
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").Value = "'0.9993"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "2.018.99"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'1.815"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "'0.5198"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("D50").Value = "'9.381"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +1.33%  "
